$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 2739
$ws.Cells.Item(2, 5).Value = 334
$ws.Cells.Item(2, 6).Value = 334
$ws.Cells.Item(2, 7).Value = 343
$ws.Cells.Item(2, 8).Value = 261
$ws.Cells.Item(2, 9).Value = 261
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 2035
$ws.Cells.Item(2, 12).Value = 311
$ws.Cells.Item(2, 13).Value = 1725
$ws.Cells.Item(2, 14).Value = 1724
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 36
$ws.Cells.Item(2, 17).Value = 222
$ws.Cells.Item(2, 18).Value = -181
$ws.Cells.Item(2, 19).Value = -47
$ws.Cells.Item(2, 20).Value = 147
$ws.Cells.Item(2, 21).Value = 75
$ws.Cells.Item(2, 22).Value = 7
$ws.Cells.Item(2, 23).Value = 12.19
$ws.Cells.Item(2, 24).Value = 9.54
$ws.Cells.Item(2, 25).Value = 16.14
$ws.Cells.Item(2, 26).Value = 13.54
$ws.Cells.Item(2, 27).Value = 18.01
$ws.Cells.Item(2, 28).Value = 4776.27
$ws.Cells.Item(2, 29).Value = 3626
$ws.Cells.Item(2, 30).Value = 9.55
$ws.Cells.Item(2, 31).Value = 24196
$ws.Cells.Item(2, 32).Value = 1.43
$ws.Cells.Item(2, 33).Value = 250
$ws.Cells.Item(2, 34).Value = 0.72
$ws.Cells.Item(2, 35).Value = 6.89
$ws.Cells.Item(2, 36).Value = 7244890

# Row 3
$ws.Cells.Item(3, 4).Value = 2704
$ws.Cells.Item(3, 5).Value = 369
$ws.Cells.Item(3, 6).Value = 369
$ws.Cells.Item(3, 7).Value = 392
$ws.Cells.Item(3, 8).Value = 306
$ws.Cells.Item(3, 9).Value = 306
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 2243
$ws.Cells.Item(3, 12).Value = 319
$ws.Cells.Item(3, 13).Value = 1924
$ws.Cells.Item(3, 14).Value = 1924
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 36
$ws.Cells.Item(3, 17).Value = 404
$ws.Cells.Item(3, 18).Value = -185
$ws.Cells.Item(3, 19).Value = -101
$ws.Cells.Item(3, 20).Value = 200
$ws.Cells.Item(3, 21).Value = 204
$ws.Cells.Item(3, 22).Value = 9
$ws.Cells.Item(3, 23).Value = 13.64
$ws.Cells.Item(3, 24).Value = 11.33
$ws.Cells.Item(3, 25).Value = 16.79
$ws.Cells.Item(3, 26).Value = 14.32
$ws.Cells.Item(3, 27).Value = 16.59
$ws.Cells.Item(3, 28).Value = 5617.64
$ws.Cells.Item(3, 29).Value = 4226
$ws.Cells.Item(3, 30).Value = 8.35
$ws.Cells.Item(3, 31).Value = 28279
$ws.Cells.Item(3, 32).Value = 1.25
$ws.Cells.Item(3, 33).Value = 250
$ws.Cells.Item(3, 34).Value = 0.71
$ws.Cells.Item(3, 35).Value = 5.55
$ws.Cells.Item(3, 36).Value = 7293670

# Row 4
$ws.Cells.Item(4, 4).Value = 951
$ws.Cells.Item(4, 5).Value = 26
$ws.Cells.Item(4, 6).Value = 405
$ws.Cells.Item(4, 7).Value = 32
$ws.Cells.Item(4, 8).Value = 345
$ws.Cells.Item(4, 9).Value = 345
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2548
$ws.Cells.Item(4, 12).Value = 383
$ws.Cells.Item(4, 13).Value = 2165
$ws.Cells.Item(4, 14).Value = 2165
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(4, 16).Value = 37
$ws.Cells.Item(4, 17).Value = 554
$ws.Cells.Item(4, 18).Value = -248
$ws.Cells.Item(4, 19).Value = -118
$ws.Cells.Item(4, 20).Value = 248
$ws.Cells.Item(4, 21).Value = 306
$ws.Cells.Item(4, 22).Value = 13
$ws.Cells.Item(4, 23).Value = 2.71
$ws.Cells.Item(4, 24).Value = 36.28
$ws.Cells.Item(4, 25).Value = 16.86
$ws.Cells.Item(4, 26).Value = 14.41
$ws.Cells.Item(4, 27).Value = 17.69
$ws.Cells.Item(4, 28).Value = 6231.44
$ws.Cells.Item(4, 29).Value = 4842
$ws.Cells.Item(4, 30).Value = 10.69
$ws.Cells.Item(4, 31).Value = 33017
$ws.Cells.Item(4, 32).Value = 1.57
$ws.Cells.Item(4, 33).Value = 550
$ws.Cells.Item(4, 34).Value = 1.06
$ws.Cells.Item(4, 35).Value = 10.46
$ws.Cells.Item(4, 36).Value = 7046410

# Row 5
$ws.Cells.Item(5, 4).Value = 1176
$ws.Cells.Item(5, 5).Value = 92
$ws.Cells.Item(5, 6).Value = 92
$ws.Cells.Item(5, 7).Value = 118
$ws.Cells.Item(5, 8).Value = 127
$ws.Cells.Item(5, 9).Value = 127
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 1638
$ws.Cells.Item(5, 12).Value = 415
$ws.Cells.Item(5, 13).Value = 1224
$ws.Cells.Item(5, 14).Value = 1223
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 13
$ws.Cells.Item(5, 17).Value = -90
$ws.Cells.Item(5, 18).Value = -246
$ws.Cells.Item(5, 19).Value = -12
$ws.Cells.Item(5, 20).Value = 122
$ws.Cells.Item(5, 21).Value = -212
$ws.Cells.Item(5, 22).Value = 50
$ws.Cells.Item(5, 23).Value = 7.8
$ws.Cells.Item(5, 24).Value = 10.83
$ws.Cells.Item(5, 25).Value = 7.5
$ws.Cells.Item(5, 26).Value = 6.08
$ws.Cells.Item(5, 27).Value = 33.88
$ws.Cells.Item(5, 28).Value = 22413.08
$ws.Cells.Item(5, 29).Value = 3331
$ws.Cells.Item(5, 30).Value = 14.56
$ws.Cells.Item(5, 31).Value = 53040
$ws.Cells.Item(5, 32).Value = 0.91
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 2331774

# Row 6
$ws.Cells.Item(6, 4).Value = 1420
$ws.Cells.Item(6, 5).Value = 36
$ws.Cells.Item(6, 6).Value = 36
$ws.Cells.Item(6, 7).Value = 132
$ws.Cells.Item(6, 8).Value = 97
$ws.Cells.Item(6, 9).Value = 97
$ws.Cells.Item(6, 11).Value = 2220
$ws.Cells.Item(6, 12).Value = 899
$ws.Cells.Item(6, 13).Value = 1321
$ws.Cells.Item(6, 14).Value = 1319
$ws.Cells.Item(6, 16).Value = 13
$ws.Cells.Item(6, 17).Value = 90
$ws.Cells.Item(6, 18).Value = -304
$ws.Cells.Item(6, 19).Value = 229
$ws.Cells.Item(6, 20).Value = 16
$ws.Cells.Item(6, 21).Value = 74
$ws.Cells.Item(6, 22).Value = 283
$ws.Cells.Item(6, 23).Value = 2.54
$ws.Cells.Item(6, 24).Value = 6.86
$ws.Cells.Item(6, 25).Value = 7.62
$ws.Cells.Item(6, 26).Value = 5.05
$ws.Cells.Item(6, 27).Value = 68.11
$ws.Cells.Item(6, 28).Value = 9983.53
$ws.Cells.Item(6, 29).Value = 4154
$ws.Cells.Item(6, 30).Value = 8.32
$ws.Cells.Item(6, 31).Value = 56573
$ws.Cells.Item(6, 32).Value = 0.61
$ws.Cells.Item(6, 33).Value = 400
$ws.Cells.Item(6, 34).Value = 1.16
$ws.Cells.Item(6, 35).Value = 9.63
$ws.Cells.Item(6, 36).Value = 2331774

# Clear rows 7-9 (D:AI), keep A/B/C
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
